$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.078.78'
$ws.Range("E2").Value = '  +0.23%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.313.01'
$ws.Range("E3").Value = '  +0.20%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.28'
$ws.Range("E5").Value = '  -0.35%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '99.09'
$ws.Range("E6").Value = '  -1.35%  '

$ws.Range("E7").Value = '  +0.99%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("E9").Value = '  +2.35%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.81'
$ws.Range("E10").Value = '  +2.38%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0790'
$ws.Range("E11").Value = '  -0.70%  '

$ws.Range("E12").Value = '  -1.03%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '17.98'
$ws.Range("E13").Value = '  +0.62%  '

$ws.Range("E14").Value = '  +0.03%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.670.10'
$ws.Range("E15").Value = '  -0.55%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.384.50'
$ws.Range("E16").Value = '  +4.26%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.791'
$ws.Range("E17").Value = '  -2.91%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.973.91'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.48'
$ws.Range("E19").Value = '  +6.73%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0913'
$ws.Range("E20").Value = '  +1.01%  '

$ws.Range("E21").Value = '  +0.75%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.94'
$ws.Range("E22").Value = '  +0.14%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '240.49'
$ws.Range("E23").Value = '  +1.43%  '

$ws.Range("E24").Value = '  -1.72%  '

$ws.Range("E25").Value = '  +0.05%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.45'
$ws.Range("E26").Value = '  -0.95%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.93'
$ws.Range("E27").Value = '  +0.73%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '169.33'
$ws.Range("E28").Value = '  +0.76%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.05'
$ws.Range("E29").Value = '  -10.43%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.20'
$ws.Range("E30").Value = '  +0.14%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '33.56'
$ws.Range("E31").Value = '  -0.99%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.24'
$ws.Range("E32").Value = '  +4.64%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.92'
$ws.Range("E33").Value = '  +6.90%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '18.36'
$ws.Range("E35").Value = '  +7.99%  '

$ws.Range("E36").Value = '  +0.13%  '

$ws.Range("E37").Value = '  +0.45%  '

$ws.Range("E38").Value = '  +0.32%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.81'
$ws.Range("E39").Value = '  +0.98%  '

$ws.Range("E40").Value = '  -2.09%  '

$ws.Range("E41").Value = '  +0.05%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.995.39'
$ws.Range("E42").Value = '  -0.32%  '

$ws.Range("E43").Value = '  +1.15%  '

$ws.Range("B44").Value = 'ApeXProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.18'
$ws.Range("E44").Value = '  -4.97%  '

$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.10'
$ws.Range("E45").Value = '  -0.87%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.56'
$ws.Range("E46").Value = '  +0.63%  '

$ws.Range("E47").Value = '  +0.29%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '54.75'
$ws.Range("E48").Value = '  -0.93%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '74.85'
$ws.Range("E49").Value = '  +7.05%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.536.26'
$ws.Range("E50").Value = '  +0.68%  '

$ws.Range("E51").Value = '  +0.27%  '
